$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new blank rows at 409-410, shifting the existing 409:432 block
# (and everything below it) down to 411:434.
$ws.Range("A409:A410").EntireRow.Insert()

# New weekly entry - "Primera" quality, Región Metropolitana
$ws.Cells.Item(409, 1).Value = 9
$ws.Cells.Item(409, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(409, 3).Value = "Metropolitana"
$ws.Cells.Item(409, 4).Value = 44516
$ws.Cells.Item(409, 5).Value = 13
$ws.Cells.Item(409, 6).Value = 100112008
$ws.Cells.Item(409, 7).Value = "Coliflor"
$ws.Cells.Item(409, 8).Value = "Sin especificar"
$ws.Cells.Item(409, 9).Value = "Primera"
$ws.Cells.Item(409, 10).Value = 4300
$ws.Cells.Item(409, 11).Value = 500
$ws.Cells.Item(409, 12).Value = 600
$ws.Cells.Item(409, 13).Value = 550
$ws.Cells.Item(409, 14).Value = "`$/unidad"
$ws.Cells.Item(409, 15).Value = "Región Metropolitana"
$ws.Cells.Item(409, 16).Value = 550
$ws.Cells.Item(409, 17).Value = 1
$ws.Cells.Item(409, 18).Value = "Hortaliza"

# New weekly entry - "Segunda" quality, Región Metropolitana
$ws.Cells.Item(410, 1).Value = 9
$ws.Cells.Item(410, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(410, 3).Value = "Metropolitana"
$ws.Cells.Item(410, 4).Value = 44516
$ws.Cells.Item(410, 5).Value = 13
$ws.Cells.Item(410, 6).Value = 100112008
$ws.Cells.Item(410, 7).Value = "Coliflor"
$ws.Cells.Item(410, 8).Value = "Sin especificar"
$ws.Cells.Item(410, 9).Value = "Segunda"
$ws.Cells.Item(410, 10).Value = 1600
$ws.Cells.Item(410, 11).Value = 400
$ws.Cells.Item(410, 12).Value = 400
$ws.Cells.Item(410, 13).Value = 400
$ws.Cells.Item(410, 14).Value = "`$/unidad"
$ws.Cells.Item(410, 15).Value = "Región Metropolitana"
$ws.Cells.Item(410, 16).Value = 400
$ws.Cells.Item(410, 17).Value = 1
$ws.Cells.Item(410, 18).Value = "Hortaliza"
